$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 222, shifting existing rows 222:237 down to 223:238
$ws.Rows("222:222").Insert()

# Populate the newly inserted row 222 with the new record
$ws.Cells.Item(222, 1).Value = 5
$ws.Cells.Item(222, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(222, 3).Value = "Maule"
$ws.Cells.Item(222, 4).Value = 45013
$ws.Cells.Item(222, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(222, 5).Value = 7
$ws.Cells.Item(222, 6).Value = 100112031
$ws.Cells.Item(222, 7).Value = "Poroto verde"
$ws.Cells.Item(222, 8).Value = "Sin especificar"
$ws.Cells.Item(222, 9).Value = "Primera"
$ws.Cells.Item(222, 10).Value = 200
$ws.Cells.Item(222, 11).Value = 30000
$ws.Cells.Item(222, 12).Value = 30000
$ws.Cells.Item(222, 13).Value = 30000
$ws.Cells.Item(222, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(222, 15).Value = "Región Metropolitana"
$ws.Cells.Item(222, 16).Value = 1200
$ws.Cells.Item(222, 17).Value = 25
$ws.Cells.Item(222, 18).Value = "Hortaliza"
